$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = 0.0148290494731055
$ws.Range("H1").Value = -0.00360507125957782
$ws.Range("I1").Value = 0.000338847568184501
$ws.Range("J1").Value = -0.0867370786746085
$ws.Range("K1").Value = "foo"

$ws.Columns.Item(5).ColumnWidth = 14.59

$ws.Range("K1").Select()
